$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rotation block (rows 3-6, columns C:D): apply a 3-decimal number format ---
$ws.Range("C3:D6").NumberFormat = "0.000_ "

# --- Translation block (rows 7-9, columns C:D): drop the *1000 helper formulas,
#     keep the original (pre-multiplication) numeric values, and apply the same
#     3-decimal number format used above ---
$ws.Range("C7").Value = 1.70079118954
$ws.Range("D7").Value = 3.412
$ws.Range("C8").Value = 0.0159456324149
$ws.Range("D8").Value = 0
$ws.Range("C9").Value = 1.51095763913
$ws.Range("D9").Value = 0.5
$ws.Range("C7:D9").NumberFormat = "0.000_ "

# --- Move the active selection to D11 ---
$null = $ws.Range("D11").Select()
